$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Cells.Item(28, 1).Value = 11.46
$ws.Cells.Item(28, 2).Value = 17.55
$ws.Cells.Item(28, 3).Value = 4.2
$ws.Cells.Item(28, 4).Value = 3.66
$ws.Cells.Item(28, 5).Value = 36.87
$ws.Cells.Item(28, 9).Value = 3.24
$ws.Cells.Item(28, 11).Value = 1.55
$ws.Cells.Item(28, 13).Value = 1.9
$ws.Cells.Item(28, 15).Value = 1.35
$ws.Cells.Item(28, 17).Value = 1.8
$ws.Cells.Item(28, 18).Value = 1.62
$ws.Cells.Item(28, 19).Value = 1.15
$ws.Cells.Item(28, 20).Value = 0.6
$ws.Cells.Item(28, 21).Value = 1.9
$ws.Cells.Item(28, 22).Value = 0.1
$ws.Cells.Item(28, 23).Value = 0.82
$ws.Cells.Item(28, 24).Value = 1.85
$ws.Cells.Item(28, 25).Value = 0.5
$ws.Cells.Item(28, 26).Value = 0.35
$ws.Cells.Item(28, 27).Value = 0.7
$ws.Cells.Item(28, 28).Value = 4.1
$ws.Cells.Item(28, 29).Value = 2.4
$ws.Cells.Item(28, 30).Value = 1.53
$ws.Cells.Item(28, 31).Value = 0.45
$ws.Cells.Item(28, 32).Value = 1.1
$ws.Cells.Item(28, 33).Value = 2
$ws.Cells.Item(28, 37).Value = 2.46
$ws.Cells.Item(28, 38).Value = 0.4
$ws.Cells.Item(28, 39).Value = 0.4
$ws.Cells.Item(28, 40).Value = 0.4
$ws.Cells.Item(28, 42).Value = 0.5
$ws.Cells.Item(28, 43).Value = 0.45
$ws.Cells.Item(28, 44).Value = 2.4

# Row 30
$ws.Cells.Item(30, 9).Value = 1.8
$ws.Cells.Item(30, 11).Value = 2.3
$ws.Cells.Item(30, 13).Value = 1.3
$ws.Cells.Item(30, 15).Value = 2.7
$ws.Cells.Item(30, 17).Value = 1.8
$ws.Cells.Item(30, 18).Value = 2.6
$ws.Cells.Item(30, 19).Value = 0.75
$ws.Cells.Item(30, 20).Value = 0.6
$ws.Cells.Item(30, 21).Value = 0.6
$ws.Cells.Item(30, 22).Value = 0.55
$ws.Cells.Item(30, 23).Value = 0.83
$ws.Cells.Item(30, 25).Value = 1.1
$ws.Cells.Item(30, 26).Value = 0.5
$ws.Cells.Item(30, 27).Value = 0.4
$ws.Cells.Item(30, 28).Value = 1.01
$ws.Cells.Item(30, 30).Value = 2.2
$ws.Cells.Item(30, 31).Value = 1
$ws.Cells.Item(30, 33).Value = 8.199999999999999
$ws.Cells.Item(30, 34).Value = 1.1
$ws.Cells.Item(30, 37).Value = 0.9
$ws.Cells.Item(30, 39).Value = 0.5
$ws.Cells.Item(30, 40).Value = 0.5
$ws.Cells.Item(30, 41).Value = 0.4
$ws.Cells.Item(30, 44).Value = 3

# Row 31
$ws.Cells.Item(31, 9).Value = 2.8
$ws.Cells.Item(31, 11).Value = 1.5
$ws.Cells.Item(31, 13).Value = 0.6
$ws.Cells.Item(31, 15).Value = 0.75
$ws.Cells.Item(31, 17).Value = 1
$ws.Cells.Item(31, 18).Value = 2.08
$ws.Cells.Item(31, 19).Value = 1.55
$ws.Cells.Item(31, 20).Value = 0.1
$ws.Cells.Item(31, 21).Value = 1.65
$ws.Cells.Item(31, 22).Value = 0.4
$ws.Cells.Item(31, 23).Value = 0.4
$ws.Cells.Item(31, 25).Value = 1.5
$ws.Cells.Item(31, 26).Value = 0.4
$ws.Cells.Item(31, 27).Value = 0.4
$ws.Cells.Item(31, 28).Value = 3.6
$ws.Cells.Item(31, 29).Value = 3.6
$ws.Cells.Item(31, 30).Value = 1
$ws.Cells.Item(31, 31).Value = 1.1
$ws.Cells.Item(31, 32).Value = 1.1
$ws.Cells.Item(31, 33).Value = 9.6
$ws.Cells.Item(31, 34).Value = 1.1
$ws.Cells.Item(31, 35).Value = 1.1
$ws.Cells.Item(31, 37).Value = 1.05
$ws.Cells.Item(31, 39).Value = 0.4
$ws.Cells.Item(31, 41).Value = 0.4
$ws.Cells.Item(31, 42).Value = 0.5
$ws.Cells.Item(31, 43).Value = 0.35

# Row 32
$ws.Cells.Item(32, 9).Value = 3.41
$ws.Cells.Item(32, 11).Value = 3.5
$ws.Cells.Item(32, 13).Value = 3.35
$ws.Cells.Item(32, 15).Value = 3.5
$ws.Cells.Item(32, 17).Value = 3.4
$ws.Cells.Item(32, 18).Value = 2.85
$ws.Cells.Item(32, 20).Value = 0.6
$ws.Cells.Item(32, 21).Value = 2.15
$ws.Cells.Item(32, 22).Value = 0.7
$ws.Cells.Item(32, 23).Value = 0.7
$ws.Cells.Item(32, 24).Value = 1.85
$ws.Cells.Item(32, 25).Value = 2.23
$ws.Cells.Item(32, 26).Value = 0.9
$ws.Cells.Item(32, 27).Value = 0.85
$ws.Cells.Item(32, 28).Value = 2.27
$ws.Cells.Item(32, 29).Value = 2.3
$ws.Cells.Item(32, 30).Value = 1.2
$ws.Cells.Item(32, 31).Value = 0.6
$ws.Cells.Item(32, 33).Value = 9.130000000000001
$ws.Cells.Item(32, 34).Value = 1.1
$ws.Cells.Item(32, 35).Value = 2.3
$ws.Cells.Item(32, 37).Value = 1
$ws.Cells.Item(32, 38).Value = 0.67
$ws.Cells.Item(32, 39).Value = 0.63
$ws.Cells.Item(32, 40).Value = 0.3
$ws.Cells.Item(32, 41).Value = 0.4
$ws.Cells.Item(32, 42).Value = 0.4
$ws.Cells.Item(32, 43).Value = 0.45
$ws.Cells.Item(32, 44).Value = 2.4

# Row 71
$ws.Cells.Item(71, 9).Value = 3.63
